# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value = "25.962.47"
$c.Style = "Normal"
$ws.Cells.Item(2, 5).Value = "  +0.23%  "
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value = "1.645.54"
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = "  +0.44%  "
$c = $ws.Cells.Item(4, 4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = "215.99"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +0.31%  "
$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = "0.5094"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +1.25%  "
$c = $ws.Cells.Item(7, 4)
$c.NumberFormat = "@"
$c.Value = "1.006"
$c.Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.33%  "
$c = $ws.Cells.Item(8, 4)
$c.NumberFormat = "@"
$c.Value = "0.2568"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  +0.13%  "
$c = $ws.Cells.Item(9, 4)
$c.NumberFormat = "@"
$c.Value = "0.06395"
$c.Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.08%  "
$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = "19.55"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -0.62%  "
$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = "0.07803"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +1.02%  "
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = "4.308"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +1.10%  "
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value = "1.648.97"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  +0.58%  "
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value = "0.5480"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +0.62%  "
$c = $ws.Cells.Item(15, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₅7869"
$c.Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  -0.47%  "
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value = "64.57"
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +0.34%  "
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value = "26.035.96"
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = "  +0.48%  "
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value = "198.35"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -2.37%  "
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = "4.452"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +1.71%  "
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = "9.992"
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  +0.86%  "
$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = "6.070"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +1.44%  "
$c = $ws.Cells.Item(23, 4)
$c.NumberFormat = "@"
$c.Value = "1.008"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +0.39%  "
$ws.Cells.Item(24, 5).Value = "  -2.38%  "
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value = "141.37"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.09%  "
$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = "0.1145"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +0.82%  "
$c = $ws.Cells.Item(27, 4)
$c.NumberFormat = "@"
$c.Value = "6.898"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.53%  "
$ws.Cells.Item(28, 5).Value = "  +0.60%  "
$ws.Cells.Item(29, 2).Value = "PancakeSwap"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = "1.243"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +0.02%  "
$ws.Cells.Item(30, 2).Value = "Hedera"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value = "0.05055"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +2.40%  "
$c = $ws.Cells.Item(31, 4)
$c.NumberFormat = "@"
$c.Value = "3.272"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -0.17%  "
$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = "3.201"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +0.55%  "
$c = $ws.Cells.Item(33, 4)
$c.NumberFormat = "@"
$c.Value = "1.548"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +0.23%  "
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = "2.371"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  -0.23%  "
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value = "0.8979"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +0.40%  "
$ws.Cells.Item(36, 5).Value = "  -1.16%  "
$c = $ws.Cells.Item(37, 4)
$c.NumberFormat = "@"
$c.Value = "1.136.50"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -2.15%  "
$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = "0.5509"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -1.71%  "
$ws.Cells.Item(39, 2).Value = "BabyDogeCoin"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$c = $ws.Cells.Item(39, 4)
$c.NumberFormat = "@"
$c.Value = "0.0₈134"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +15.53%  "
$ws.Cells.Item(40, 2).Value = "VeChain"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = "0.01562"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.25%  "
$c = $ws.Cells.Item(41, 4)
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.44%  "
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = "2.549"
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.55%  "
$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = "5.636"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -1.45%  "
$c = $ws.Cells.Item(44, 4)
$c.NumberFormat = "@"
$c.Value = "0.8224"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +1.80%  "
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = "100.28"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.46%  "
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = "1.781.20"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +0.27%  "
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = "0.4543"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.40%  "
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value = "1.004"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.22%  "
$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = "55.13"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +0.22%  "
$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = "0.05081"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.53%  "
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = "1.007"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.40%  "

Write-Output "Applied 103 cell updates"
